$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 47098
$ws.Range("J3").Value = 47098
$ws.Range("L3").Value = 47098
$ws.Range("N3").Value = -47326
$ws.Range("H39").Value = 150.10527
$ws.Range("I39").Value = 80.85714
$ws.Range("J39").Value = 190.5
$ws.Range("K39").Value = 242.57142
$ws.Range("L39").Value = 571.5
$ws.Range("M39").Value = 53.42858000000001
$ws.Range("N39").Value = -1163.5
$ws.Range("H62").Value = 1921.6471
$ws.Range("I62").Value = 2004.8572
$ws.Range("J62").Value = 1533.3334
$ws.Range("K62").Value = 2004.8572
$ws.Range("L62").Value = 1533.3334
$ws.Range("M62").Value = -1380.8572
$ws.Range("N62").Value = -2781.3334
$ws.Range("H65").Value = 1921.6471
$ws.Range("I65").Value = 2004.8572
$ws.Range("J65").Value = 1533.3334
$ws.Range("K65").Value = 10024.286
$ws.Range("L65").Value = 7666.666999999999
$ws.Range("M65").Value = -6904.286
$ws.Range("N65").Value = -13906.667
$ws.Range("H92").Value = 4247.8887
$ws.Range("I92").Value = 1397.4667
$ws.Range("J92").Value = 18500
$ws.Range("K92").Value = 1397.4667
$ws.Range("L92").Value = 18500
$ws.Range("M92").Value = -149.4666999999999
$ws.Range("N92").Value = -20996
$ws.Range("H102").Value = 47098
$ws.Range("J102").Value = 47098
$ws.Range("L102").Value = 47098
$ws.Range("N102").Value = -53588
$ws.Range("H103").Value = 7513713
$ws.Range("I103").Value = 10017701
$ws.Range("J103").Value = 1750
$ws.Range("K103").Value = 30053103
$ws.Range("L103").Value = 5250
$ws.Range("M103").Value = -30052517
$ws.Range("N103").Value = -6422
$ws.Range("H107").Value = 6370.375
$ws.Range("I107").Value = 6790.4
$ws.Range("J107").Value = 70
$ws.Range("K107").Value = 6790.4
$ws.Range("L107").Value = 70
$ws.Range("M107").Value = -4870.4
$ws.Range("N107").Value = -3910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9805787
$ws.Range("I61").Value = 11112861
$ws.Range("J61").Value = 2733.3333
$ws.Range("K61").Value = 11112861
$ws.Range("L61").Value = 2733.3333
$ws.Range("M61").Value = -11112649
$ws.Range("N61").Value = -3157.3333
$ws.Range("H74").Value = 13890810
$ws.Range("I74").Value = 17242514
$ws.Range("K74").Value = 17242514
$ws.Range("M74").Value = -17241640
$ws.Range("H77").Value = 13890810
$ws.Range("I77").Value = 17242514
$ws.Range("K77").Value = 86212570
$ws.Range("M77").Value = -86208202
$ws.Range("H102").Value = 1333.3334
$ws.Range("I102").Value = 1433.3334
$ws.Range("J102").Value = 1233.3334
$ws.Range("K102").Value = 1433.3334
$ws.Range("L102").Value = 1233.3334
$ws.Range("M102").Value = 188.6666
$ws.Range("N102").Value = -4477.3334
$ws.Range("H136").Value = 9805787
$ws.Range("I136").Value = 11112861
$ws.Range("J136").Value = 2733.3333
$ws.Range("K136").Value = 33338583
$ws.Range("L136").Value = 8199.999899999999
$ws.Range("M136").Value = -33336033
$ws.Range("N136").Value = -13299.9999
$ws.Range("H139").Value = 56531.668
$ws.Range("J139").Value = 56531.668
$ws.Range("L139").Value = 56531.668
$ws.Range("N139").Value = -66811.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 47636.668
$ws.Range("J53").Value = 47636.668
$ws.Range("L53").Value = 47636.668
$ws.Range("N53").Value = -48784.668
$ws.Range("H81").Value = 20855.2
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 20855.2
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20855.2
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -22977.2
$ws.Range("H84").Value = 20855.2
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 20855.2
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 62565.60000000001
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -73173.60000000001
$ws.Range("H118").Value = 7846.6665
$ws.Range("J118").Value = 7846.6665
$ws.Range("L118").Value = 7846.6665
$ws.Range("N118").Value = -11160.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 24500
$ws.Range("I45").Value = 15000
$ws.Range("J45").Value = 29250
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 29250
$ws.Range("M45").Value = -14407
$ws.Range("N45").Value = -30436
$ws.Range("H108").Value = 20220.75
$ws.Range("J108").Value = 20220.75
$ws.Range("L108").Value = 20220.75
$ws.Range("N108").Value = -27900.75
$ws.Range("H127").Value = 51994.445
$ws.Range("J127").Value = 51994.445
$ws.Range("L127").Value = 51994.445
$ws.Range("N127").Value = -61914.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 1284.7548
$ws.Range("J88").Value = 1284.7548
$ws.Range("L88").Value = 3854.2644
$ws.Range("N88").Value = -4710.2644
$ws.Range("H91").Value = 1284.7548
$ws.Range("J91").Value = 1284.7548
$ws.Range("L91").Value = 3854.2644
$ws.Range("N91").Value = -6818.2644
$ws.Range("H92").Value = 760.8
$ws.Range("I92").Value = 634.6667
$ws.Range("J92").Value = 950
$ws.Range("K92").Value = 1904.0001
$ws.Range("L92").Value = 2850
$ws.Range("M92").Value = -656.0001
$ws.Range("N92").Value = -5346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1115.5416
$ws.Range("I97").Value = 957.5
$ws.Range("J97").Value = 1431.625
$ws.Range("K97").Value = 957.5
$ws.Range("L97").Value = 1431.625
$ws.Range("M97").Value = -461.5
$ws.Range("N97").Value = -2423.625
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1369.1052
$ws.Range("I93").Value = 1421.9166
$ws.Range("J93").Value = 1278.5714
$ws.Range("K93").Value = 1421.9166
$ws.Range("L93").Value = 1278.5714
$ws.Range("M93").Value = -173.9166
$ws.Range("N93").Value = -3774.5714
$ws.Range("H100").Value = 3038.182
$ws.Range("I100").Value = 2995
$ws.Range("J100").Value = 3090
$ws.Range("K100").Value = 2995
$ws.Range("L100").Value = 3090
$ws.Range("M100").Value = -2454
$ws.Range("N100").Value = -4172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 831.8946999999999
$ws.Range("I107").Value = 1320.3334
$ws.Range("J107").Value = 392.3
$ws.Range("K107").Value = 3961.0002
$ws.Range("L107").Value = 1176.9
$ws.Range("M107").Value = -2041.0002
$ws.Range("N107").Value = -5016.9
